$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay as text (matches original inlineStr type).
# Force text format, set value, then reset style so no stray formatting remains.
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D16","D17","D18","D19","D21","D22","D24","D25","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.417.73"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.794.55"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "338.26"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "0.3797"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").Value = "0.3457"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").Value = "48.88"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "1.199"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "0.07513"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "21.99"
$ws.Range("E13").Value = "  +6.88%  "
$ws.Range("D14").Value = "6.468"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "1.794.77"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "7.075"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "0.00001103"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "0.06668"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "84.91"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "6.528"
$ws.Range("E21").Value = "  +4.26%  "
$ws.Range("D22").Value = "17.37"
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("D23").Value = "27.419.00"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "12.51"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").Value = "2.433"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "2.562"
$ws.Range("E26").Value = "  +5.08%  "
$ws.Range("D27").Value = "1.486"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "21.47"
$ws.Range("E28").Value = "  +8.81%  "
$ws.Range("D29").Value = "153.70"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "1.998.21"
$ws.Range("D31").Value = "133.62"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "4.068"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").Value = "6.090"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "0.08705"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "13.25"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").Value = "1.655"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "5.457"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.6904"
$ws.Range("E38").Value = "  +8.50%  "
$ws.Range("D39").Value = "0.06388"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").Value = "8.886"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").Value = "0.02348"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("E43").Value = "  +3.81%  "
$ws.Range("D44").Value = "14.43"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "0.6439"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").Value = "3.873"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "2.135"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").Value = "130.02"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "0.07200"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "79.57"
$ws.Range("E51").Value = "  +1.31%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
